# "Se procesan de nuevo los datos con las nuevas dimensiones curadas"
#
# The municipio-nombre and aragon columns are re-processed against the
# newly curated dimensions: both now point at the generic
# sdmx-dimension:refArea dimension (instead of their bespoke
# iaest-measure:/iaest-dimension: values), their "medida" row collapses to
# "dim" to match, the aragon column's concept-scheme marker becomes the new
# URI-Comunidad (instead of skos:Concept), and the now-obsolete
# mapping-aragon.xlsx mapping file reference is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# municipio-nombre (column D) re-measured as a standard reference-area dimension
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"

# aragon (column F) re-measured as a standard reference-area dimension
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F4").Value = "URI-Comunidad"

# the old per-dataset mapping-aragon.xlsx reference no longer applies, so
# remove the cell outright (not just blank its text)
$ws.Range("F5").Clear()
